$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2306.2666
$ws.Range("I40").Value = 3700
$ws.Range("J40").Value = 1799.4546
$ws.Range("K40").Value = 3700
$ws.Range("L40").Value = 1799.4546
$ws.Range("M40").Value = -3525
$ws.Range("N40").Value = -2149.4546
$ws.Range("H62").Value = 9500
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("H65").Value = 9500
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("H113").Value = 3051.3635
$ws.Range("I113").Value = 2610.8333
$ws.Range("J113").Value = 3580
$ws.Range("K113").Value = 2610.8333
$ws.Range("L113").Value = 3580
$ws.Range("M113").Value = 643.1667000000002
$ws.Range("N113").Value = -10088
$ws.Range("H137").Value = 1166.25
$ws.Range("I137").Value = 939.1429000000001
$ws.Range("J137").Value = 1847.5714
$ws.Range("K137").Value = 2817.4287
$ws.Range("L137").Value = 5542.7142
$ws.Range("M137").Value = -267.4287000000004
$ws.Range("N137").Value = -10642.7142
$ws.Range("H138").Value = 474236.16
$ws.Range("I138").Value = 1497.826
$ws.Range("K138").Value = 4493.478
$ws.Range("M138").Value = 646.5219999999999
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7046.8667
$ws.Range("I2").Value = 299.23077
$ws.Range("K2").Value = 299.23077
$ws.Range("M2").Value = -186.23077
$ws.Range("H61").Value = 811.17145
$ws.Range("I61").Value = 599.6896400000001
$ws.Range("K61").Value = 599.6896400000001
$ws.Range("M61").Value = -387.6896400000001
$ws.Range("H74").Value = 1347.36
$ws.Range("I74").Value = 637.26666
$ws.Range("J74").Value = 2412.5
$ws.Range("K74").Value = 637.26666
$ws.Range("L74").Value = 2412.5
$ws.Range("M74").Value = 236.73334
$ws.Range("N74").Value = -4160.5
$ws.Range("H77").Value = 1347.36
$ws.Range("I77").Value = 637.26666
$ws.Range("J77").Value = 2412.5
$ws.Range("K77").Value = 3186.3333
$ws.Range("L77").Value = 12062.5
$ws.Range("M77").Value = 1181.6667
$ws.Range("N77").Value = -20798.5
$ws.Range("H116").Value = 7046.8667
$ws.Range("I116").Value = 299.23077
$ws.Range("K116").Value = 299.23077
$ws.Range("M116").Value = 1994.76923
$ws.Range("H132").Value = 1746.2826
$ws.Range("I132").Value = 1467.0857
$ws.Range("K132").Value = 4401.257100000001
$ws.Range("M132").Value = -1871.257100000001
$ws.Range("H136").Value = 811.17145
$ws.Range("I136").Value = 599.6896400000001
$ws.Range("K136").Value = 1799.06892
$ws.Range("M136").Value = 750.9310799999998
$ws.Range("H139").Value = 35192.223
$ws.Range("J139").Value = 35192.223
$ws.Range("L139").Value = 35192.223
$ws.Range("N139").Value = -45472.223

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7046.8667
$ws.Range("I3").Value = 299.23077
$ws.Range("K3").Value = 299.23077
$ws.Range("M3").Value = -185.23077
$ws.Range("H80").Value = 509.48
$ws.Range("J80").Value = 654
$ws.Range("L80").Value = 654
$ws.Range("N80").Value = -2650
$ws.Range("H83").Value = 509.48
$ws.Range("J83").Value = 654
$ws.Range("L83").Value = 3270
$ws.Range("N83").Value = -13254
$ws.Range("H132").Value = 41899.7
$ws.Range("J132").Value = 41899.7
$ws.Range("L132").Value = 41899.7
$ws.Range("N132").Value = -52019.7
$ws.Range("H133").Value = 60780
$ws.Range("J133").Value = 60780
$ws.Range("L133").Value = 60780
$ws.Range("N133").Value = -70900
$ws.Range("H134").Value = 4829.8438
$ws.Range("I134").Value = 872.4815
$ws.Range("J134").Value = 26199.6
$ws.Range("K134").Value = 2617.4445
$ws.Range("L134").Value = 78598.79999999999
$ws.Range("M134").Value = -82.44450000000006
$ws.Range("N134").Value = -83668.79999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 212.8
$ws.Range("J7").Value = 348
$ws.Range("L7").Value = 348
$ws.Range("N7").Value = -574
$ws.Range("H31").Value = 2181.5
$ws.Range("I31").Value = 1997.1428
$ws.Range("K31").Value = 1997.1428
$ws.Range("M31").Value = -1702.1428
$ws.Range("H34").Value = 2181.5
$ws.Range("I34").Value = 1997.1428
$ws.Range("K34").Value = 1997.1428
$ws.Range("M34").Value = -1795.1428
$ws.Range("H58").Value = 1302.4667
$ws.Range("I58").Value = 1003.1111
$ws.Range("K58").Value = 1003.1111
$ws.Range("M58").Value = -800.1111
$ws.Range("H86").Value = 3937207.8
$ws.Range("I86").Value = 8336922
$ws.Range("J86").Value = 26350.889
$ws.Range("K86").Value = 8336922
$ws.Range("L86").Value = 26350.889
$ws.Range("M86").Value = -8335799
$ws.Range("N86").Value = -28596.889
$ws.Range("H89").Value = 3937207.8
$ws.Range("I89").Value = 8336922
$ws.Range("J89").Value = 26350.889
$ws.Range("K89").Value = 41684610
$ws.Range("L89").Value = 131754.445
$ws.Range("M89").Value = -41678994
$ws.Range("N89").Value = -142986.445
$ws.Range("H122").Value = 813.7368
$ws.Range("J122").Value = 808
$ws.Range("L122").Value = 2424
$ws.Range("N122").Value = -7324
$ws.Range("H136").Value = 1302.4667
$ws.Range("I136").Value = 1003.1111
$ws.Range("K136").Value = 3009.3333
$ws.Range("M136").Value = -459.3332999999998
$ws.Range("H138").Value = 189635
$ws.Range("J138").Value = 189635
$ws.Range("L138").Value = 189635
$ws.Range("N138").Value = -199915
$ws.Range("H140").Value = 39900
$ws.Range("J140").Value = 39900
$ws.Range("L140").Value = 39900
$ws.Range("N140").Value = -50260

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 2028
$ws.Range("J130").Value = 2028
$ws.Range("L130").Value = 6084
$ws.Range("N130").Value = -16124
$ws.Range("H131").Value = 19609196
$ws.Range("J131").Value = 1473.7174
$ws.Range("L131").Value = 4421.1522
$ws.Range("N131").Value = -14501.1522
$ws.Range("H132").Value = 996.3333
$ws.Range("J132").Value = 1018
$ws.Range("L132").Value = 9162
$ws.Range("N132").Value = -14222
$ws.Range("H136").Value = 2396.2727
$ws.Range("I136").Value = 1407.5
$ws.Range("K136").Value = 4222.5
$ws.Range("M136").Value = 877.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1779.4375
$ws.Range("I102").Value = 1734.75
$ws.Range("K102").Value = 1734.75
$ws.Range("M102").Value = -112.75
$ws.Range("H126").Value = 1717.5555
$ws.Range("I126").Value = 1431.1111
$ws.Range("K126").Value = 4293.3333
$ws.Range("M126").Value = -1823.3333
$ws.Range("H132").Value = 2050.9
$ws.Range("I132").Value = 1846.909
$ws.Range("J132").Value = 2611.875
$ws.Range("K132").Value = 5540.727000000001
$ws.Range("L132").Value = 7835.625
$ws.Range("M132").Value = -3010.727000000001
$ws.Range("N132").Value = -12895.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2520.8
$ws.Range("I7").Value = 2366.3333
$ws.Range("J7").Value = 2752.5
$ws.Range("K7").Value = 2366.3333
$ws.Range("L7").Value = 2752.5
$ws.Range("M7").Value = -2254.3333
$ws.Range("N7").Value = -2976.5
$ws.Range("H16").Value = 709.4
$ws.Range("I16").Value = 712.1111
$ws.Range("K16").Value = 712.1111
$ws.Range("M16").Value = -542.1111
$ws.Range("H40").Value = 2692.111
$ws.Range("I40").Value = 2375.5625
$ws.Range("K40").Value = 2375.5625
$ws.Range("M40").Value = -2239.5625
$ws.Range("H93").Value = 918.4
$ws.Range("I93").Value = 622.08
$ws.Range("K93").Value = 622.08
$ws.Range("M93").Value = 625.92
$ws.Range("H126").Value = 2520.8
$ws.Range("I126").Value = 2366.3333
$ws.Range("J126").Value = 2752.5
$ws.Range("K126").Value = 7098.999899999999
$ws.Range("L126").Value = 8257.5
$ws.Range("M126").Value = -4628.999899999999
$ws.Range("N126").Value = -13197.5
$ws.Range("H136").Value = 3927.9033
$ws.Range("I136").Value = 4361.6665
$ws.Range("K136").Value = 13084.9995
$ws.Range("M136").Value = -10534.9995
$ws.Range("H137").Value = 34695.332
$ws.Range("J137").Value = 34695.332
$ws.Range("L137").Value = 34695.332
$ws.Range("N137").Value = -44895.332

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 618.4
$ws.Range("I136").Value = 567.38464
$ws.Range("K136").Value = 1702.15392
$ws.Range("M136").Value = 847.84608
